# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    (this text is shared across Overview/zh-cn/de-de since they reuse the same string)
#  - zh-cn and de-de sheets get their "Latest Target File" / "Latest Handback File" /
#    "Latest Handback DateTime" columns populated, with a hyperlink added on the
#    "Latest Target File" cell
#  - Several columns are widened to fit the newly-populated long file names

$wb = $excel.ActiveWorkbook

$overviewWs = $wb.Worksheets.Item("Overview")
$zhWs = $wb.Worksheets.Item("zh-cn")
$deWs = $wb.Worksheets.Item("de-de")

# --- Status text update (shared string used by Overview!E2/F2 and the
#     zh-cn/de-de "Status" column C2) ---
$overviewWs.Range("E2").Value = "Handed back: in sync with en-US"
$overviewWs.Range("F2").Value = "Handed back: in sync with en-US"
$zhWs.Range("C2").Value = "Handed back: in sync with en-US"
$deWs.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: populate handback columns on row 2 ---
$zhWs.Range("J2").Value = "ae7f0526-159b-4eaf-aafd-6e77a2be2935.md"
$zhWs.Hyperlinks.Add($zhWs.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0729dc005dfb2c635e2cf1a74b23e5cacd7ace06/e2e/ae7f0526-159b-4eaf-aafd-6e77a2be2935.md", "", "", "ae7f0526-159b-4eaf-aafd-6e77a2be2935.md") | Out-Null
$zhWs.Range("K2").Value = "ae7f0526-159b-4eaf-aafd-6e77a2be2935.812c5817f533bb785c302ff8fe050c1d4d1fd7d2.zh-cn.xlf"
$zhWs.Range("L2").Value = "2017-02-09 14:16:16"

# --- de-de sheet: populate handback columns on row 2 ---
$deWs.Range("J2").Value = "ae7f0526-159b-4eaf-aafd-6e77a2be2935.md"
$deWs.Hyperlinks.Add($deWs.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0729dc005dfb2c635e2cf1a74b23e5cacd7ace06/e2e/ae7f0526-159b-4eaf-aafd-6e77a2be2935.md", "", "", "ae7f0526-159b-4eaf-aafd-6e77a2be2935.md") | Out-Null
$deWs.Range("K2").Value = "ae7f0526-159b-4eaf-aafd-6e77a2be2935.812c5817f533bb785c302ff8fe050c1d4d1fd7d2.de-de.xlf"
$deWs.Range("L2").Value = "2017-02-09 14:16:42"

# --- Column width adjustments (to fit the newly-written long strings) ---
$overviewWs.Columns.Item(5).ColumnWidth = 29.9777050018311
$overviewWs.Columns.Item(6).ColumnWidth = 29.9777050018311

$zhWs.Columns.Item(3).ColumnWidth = 29.9777050018311
$zhWs.Columns.Item(10).ColumnWidth = 40
$zhWs.Columns.Item(11).ColumnWidth = 40

$deWs.Columns.Item(3).ColumnWidth = 29.9777050018311
$deWs.Columns.Item(10).ColumnWidth = 40
$deWs.Columns.Item(11).ColumnWidth = 40
